$d = $word.ActiveDocument

function Set-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# Polish diacritic helper characters (kept as codepoints to survive any
# console/encoding round-trip issues; the actual saved XML is UTF-8 correct)
$c_c_acute = [char]0x0107   # c with acute - ć
$c_z_dot   = [char]0x017C   # z with dot above - ż
$c_e_ogon  = [char]0x0119   # e with ogonek - ę

# 1. "iminfo" paragraph -> wrap word in proofErr spellStart/spellEnd
$xml1 = '<w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>iminfo</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body>'
Set-ParaXml 4 $xml1

# 2. "Nauczyc sie uzywac freqz, filter, chebord1 chebord2 yulewalk itd." -> split runs with proofErr
$xml2 = '<w:body><w:p>' +
  '<w:r><w:t xml:space="preserve">Nauczy' + $c_c_acute + ' si' + $c_e_ogon + ' </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">u' + $c_z_dot + 'ywa' + $c_c_acute + ' </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>freqz</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>filter</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve">, chebord1 chebord2 </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>yulewalk</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> itd.</w:t></w:r>' +
  '</w:p></w:body>'
Set-ParaXml 17 $xml2

# 3. "Nauczyc sie jeszcze rysowac zera i bieguny, fft itd." -> split with proofErr around fft
$xml3 = '<w:body><w:p>' +
  '<w:r><w:t xml:space="preserve">Nauczy' + $c_c_acute + ' si' + $c_e_ogon + ' jeszcze rysowa' + $c_c_acute + ' zera i bieguny, </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>fft</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> itd.</w:t></w:r>' +
  '</w:p></w:body>'
Set-ParaXml 18 $xml3

# 4. "Forma 1 sprawko zbindowane i plytka" -> split with proofErr around plytka
$xml4 = '<w:body><w:p>' +
  '<w:r><w:t xml:space="preserve">Forma 1 sprawko zbindowane i </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>plytka</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '</w:p></w:body>'
Set-ParaXml 21 $xml4

# 5. "Opisywac co to jest filtr, jaki to filtr, itd." -> split with proofErr around Opisywac
#    plus append new paragraphs after it.
$xml5 = '<w:body>' +
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Opisywac</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> co to jest filtr, jaki to filtr, itd.</w:t></w:r></w:p>' +
  '<w:p/>' +
  '<w:p><w:r><w:t>Do sprawka 9:00, 25.06, kolos 29.06</w:t></w:r></w:p>' +
  '<w:p/>' +
  '<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Rozpoznawanie twarzy </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dokonczyc</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Dodac</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> kody pod </w:t></w:r><w:r><w:t>zdj' + $c_e_ogon + 'cia</w:t></w:r></w:p>' +
  '<w:p><w:r><w:t xml:space="preserve">Rozpoznawanie </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gradientow</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
  '<w:p/>' +
  '<w:p/>' +
  '<w:p/>' +
  '</w:body>'
Set-ParaXml 24 $xml5

Write-Host "DONE ALL"
